# plans_subjects.xlsx bulk-upload cleanup
# ------------------------------------------------------------------
# The "plan_semester" column (D) had been populated with running
# semester numbers (1..10) per plan/level instead of the intended
# "semester within the level" value (1 for the first/odd semester of
# a level, 2 for the second/even one). It also had a handful of rows
# (71-92) stored as the literal text "7,8,9,10" instead of a proper
# numeric semester value. This script corrects column D to the right
# numeric values, which as a side effect drops the now-unused
# "7,8,9,10" shared string, and also freezes the header row so the
# column headers stay visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (row, corrected plan_semester value) for column D.
$changes = @(
    @(14,2),
    @(16,1), @(17,1), @(18,1), @(19,1), @(20,1), @(21,1),
    @(23,2), @(24,2), @(25,2), @(26,2), @(27,2), @(28,2), @(29,2),
    @(31,1), @(32,1), @(33,1), @(34,1), @(35,1), @(36,1),
    @(38,2), @(39,2), @(40,2), @(41,2), @(42,2), @(43,2),
    @(45,1), @(46,1), @(47,1), @(48,1), @(49,1), @(50,1),
    @(52,2), @(53,2), @(54,2), @(55,2), @(56,2), @(57,2),
    @(59,1), @(60,1), @(61,1), @(62,1), @(63,1), @(64,1),
    @(66,2), @(67,2), @(68,2), @(69,2),
    @(71,1), @(72,1), @(73,1), @(74,1), @(75,1), @(76,1), @(77,1),
    @(78,2), @(79,2), @(80,2), @(81,2), @(82,2), @(83,2), @(84,2),
    @(85,1), @(86,1),
    @(87,2), @(88,2), @(89,2),
    @(90,1), @(91,1), @(92,1),
    @(113,1), @(114,1), @(115,1), @(116,1), @(117,1), @(118,1), @(119,1), @(120,1), @(121,1), @(122,1),
    @(124,2), @(125,2), @(126,2), @(127,2), @(128,2), @(129,2), @(130,2), @(131,2), @(132,2), @(133,2), @(134,2),
    @(153,1), @(154,1), @(155,1), @(156,1), @(157,1), @(158,1), @(159,1),
    @(161,2), @(162,2), @(163,2), @(164,2), @(165,2), @(166,2),
    @(185,1), @(186,1), @(187,1), @(188,1), @(189,1), @(190,1), @(191,1), @(192,1), @(193,1),
    @(195,2), @(196,2), @(197,2), @(198,2), @(199,2), @(200,2), @(201,2), @(202,2), @(203,2)
)

foreach ($pair in $changes) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 4).Value = $val
}

# Freeze the header row (row 1) and leave the active cell on D1, as in
# the saved workbook: top pane is row 1, bottom (scrollable) pane
# starts at A2, with the selection parked on D1 in that bottom pane.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D1").Select()
